# Generate Report for Handback
# Mark the "ea844430-..." file as handed back (status flips from
# "Ready for handoff" to "Handed back: in sync with en-US") on the
# Overview sheet and on each locale sheet, and refresh the
# "Latest Handback DateTime" for both files on each locale sheet.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $handedBack
$zhcn.Range("H2").Value = "2016-03-24 03:13:14"
$zhcn.Range("H3").Value = "2016-03-24 03:13:14"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $handedBack
$dede.Range("H2").Value = "2016-03-24 03:13:28"
$dede.Range("H3").Value = "2016-03-24 03:13:28"
